$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.153.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.49%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.585.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.67%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '204.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.87%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '567.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.77%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.580.53'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.57%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.614'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.680'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.24%  '

$ws.Range('E11').Value = '  +8.03%  '

$ws.Range('E12').Value = '  -1.28%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000279'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.72%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.156.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.67%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.596.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.77%  '

$ws.Range('E17').Value = '  +1.00%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.85'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.25%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '67.938.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.12%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.82%  '

$ws.Range('E21').Value = '  +0.13%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '401.58'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.27%  '

$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.33%  '

$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.20%  '

$ws.Range('E25').Value = '  -1.14%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.43%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.33%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.90%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.33%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.26%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.44'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.05%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '666.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.40%  '

$ws.Range('E33').Value = '  -0.75%  '

$ws.Range('E34').Value = '  -0.93%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '63.25'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.05%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '41.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.69%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.409'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.63%  '

$ws.Range('E38').Value = '  -0.16%  '

$ws.Range('E39').Value = '  +11.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0750'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.57%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.175.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.09%  '

$ws.Range('E42').Value = '  +0.05%  '

$ws.Range('E43').Value = '  -0.11%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.66%  '

$ws.Range('E45').Value = '  +11.87%  '

$ws.Range('E46').Value = '  -0.32%  '

$ws.Range('E47').Value = '  +0.26%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.70'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.38%  '

$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.60'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +10.22%  '

$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.05'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.23%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '139.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.16%  '
